$d = $word.ActiveDocument

# 1. Apply double line spacing (w:spacing w:line="480" w:lineRule="auto")
#    to every paragraph in the document. LineSpacing is expressed in
#    points; Word stores twips (points * 20), so 24pt -> 480 twips.
foreach ($p in $d.Paragraphs) {
    $p.LineSpacingRule = 5   # wdLineSpaceExactly-style "auto" multiple rule
    $p.LineSpacing = 24
}

# 2. Remove the stray "_GoBack" bookmark that sits at the very start of
#    the document (first paragraph).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 3. Re-insert the "_GoBack" bookmark further down, right after "...its
#    percen" or before "tage of success..." in the second paragraph —
#    this naturally splits the run in two, matching the target markup.
$rng = $d.Content
$found = $rng.Find.Execute("Although we didn" + [char]0x2019 + "t look at live projects to see its percen", `
                            $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng)
